# Commit message: "rename file, version should not be added to filename,
# neither space or any special characters"
#
# The genome assembly identifier stored on the Target sheet carried a
# patch/version suffix (".p7"). Per the rename convention described in the
# commit, the version suffix is stripped so downstream file names built
# from this value stay clean.

$wb = $excel.ActiveWorkbook

# --- Target sheet: fix the assembly identifier -------------------------
$target = $wb.Worksheets.Item("Target")
$target.Range("D2").Value = "GRCh38"

# --- ExperimentLayout sheet: cell_pool column was storing a combined
# "guide.pool" number (e.g. 69.1, 69.3, 73.2, 73.3); it is corrected to
# hold just the pool number on its own. -----------------------------
$layout = $wb.Worksheets.Item("ExperimentLayout")

$layout.Range("E3:E21").Value = 1
$layout.Range("E22:E37").Value = 3
$layout.Range("E38:E73").Value = 2
$layout.Range("E74:E81").Value = 3
$layout.Range("E90:E96").Value = 3
$layout.Range("E99:E121").Value = 3
